# LOB1013.xlsx rebuild: the "Ementa" sheet's Objetivos/Programa
# resumido/Programa/Avaliação block got scrambled by the source CMS
# export (course code/name field shifted up a few rows, the long
# bibliography + programme text dropped, and the trailing "Requisitos"
# value row removed). Reproduce the resulting cell layout exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old trailing row 24 (the "Requisitos" value row) disappears
# entirely and the sheet's used range shrinks from A1:C24 to A1:C23 -
# drop the row outright so the dimension + row count line up.
$ws.Rows.Item(24).Delete()

# --- Row 10: Objetivos value overwritten with the docente code/name ---
$ws.Range("B10").Value = "6376612 - Daisy Rafaela da Silva"
$ws.Range("C10").Value = "6376612 - Daisy Rafaela da Silva"

# --- Row 13: was a blank-label row holding the old docente string;
# now becomes "Programa resumido:" / "Semestral" ---
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# --- Row 14: label becomes "Short syllabus:", old long description
# text is dropped (cell removed outright, not just blanked) ---
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Clear()
$ws.Range("C14").Clear()

# --- Row 15: label becomes "Programa:", value becomes the activation
# date string - must stay TEXT, not get auto-parsed into a date value ---
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2012"
$ws.Range("B15").NumberFormat = "General"
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "01/01/2012"
$ws.Range("C15").NumberFormat = "General"
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4122)

# --- Row 16: label becomes "Syllabus:", the long "Unidade primeira..."
# programme text is dropped entirely (cell removed outright) ---
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Clear()
$ws.Range("C16").Clear()

# --- Row 17: label becomes "Avaliação:", no B/C value; row reverts to
# default height ---
$ws.Range("A17").Value = "Avaliação:"
$ws.Rows.Item(17).AutoFit()

# --- Row 18: label becomes "Método:", value becomes the docente string
# again (duplicated, matching the scrambled source). B18/C18 are brand
# new cells, so copy column formats across first, then write the text. ---
$ws.Range("B19").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "6376612 - Daisy Rafaela da Silva"
$ws.Range("C18").Value = "6376612 - Daisy Rafaela da Silva"

# --- Row 19: label becomes "Critério:" (B/C keep the existing
# "A média semestral..." text - unchanged) ---
$ws.Range("A19").Value = "Critério:"

# --- Row 20: label becomes "Norma de recuperação:" (B/C keep
# "(PS+T) / 2" - unchanged) ---
$ws.Range("A20").Value = "Norma de recuperação:"

# --- Row 21: label becomes "Bibliografia:" (B/C keep the
# "Trabalho escrito..." text - unchanged) ---
$ws.Range("A21").Value = "Bibliografia:"

# --- Row 22: label becomes "Requisitos:", old CAPRA bibliography text
# dropped (cell removed outright); row reverts to default height ---
$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").Clear()
$ws.Range("C22").Clear()
$ws.Rows.Item(22).AutoFit()

# --- Row 23: old "Requisitos:" label removed (cell cleared outright),
# B/C now hold the "LOB1018 - Física I (Requisito)" text that used to
# live on row 24. B23/C23 are brand new cells here too. ---
$ws.Range("A23").Clear()
$ws.Range("B19").Copy()
$ws.Range("B23").PasteSpecial(-4122)
$ws.Range("C19").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("B23").Value = "LOB1018 -  Física I  (Requisito)`n"
$ws.Range("C23").Value = "LOB1018 -  Física I  (Requisito)`n"

# --- Row height touch-ups to match the new layout ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(23).RowHeight = 30

$ws.Range("A1").Select()
